$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "AK_1.png"
$ws.Range("B1").Value = 0.455
$ws.Range("C1").Value = 0.002
$ws.Range("D1").Value = 0.88
$ws.Range("E1").Value = 0.714
$ws.Range("F1").Value = "Akhlak Kamiswara"
$ws.Range("G1").Value = "Benar"

# Row 2
$ws.Range("A2").Value = "AK_2.png"
$ws.Range("B2").Value = 0.866
$ws.Range("C2").Value = 0.003
$ws.Range("D2").Value = 0.888
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "Akhlak Kamiswara"
$ws.Range("G2").Value = "Benar"

# Row 3
$ws.Range("A3").Value = "AK_3.png"
$ws.Range("B3").Value = 0.838
$ws.Range("C3").Value = 0.003
$ws.Range("D3").Value = 0.904
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Akhlak Kamiswara"
$ws.Range("G3").Value = "Benar"

# Row 4
$ws.Range("A4").Value = "AK_4.png"
$ws.Range("B4").Value = 0.434
$ws.Range("C4").Value = 0.001
$ws.Range("D4").Value = 0.89
$ws.Range("E4").Value = 0.857
$ws.Range("F4").Value = "Akhlak Kamiswara"
$ws.Range("G4").Value = "Benar"

# Row 5
$ws.Range("A5").Value = "AK_5.png"
$ws.Range("B5").Value = 0.441
$ws.Range("C5").Value = 0.001
$ws.Range("D5").Value = 0.899
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "Akhlak Kamiswara"
$ws.Range("G5").Value = "Benar"

# Row 6
$ws.Range("A6").Value = "MIB_1.png"
$ws.Range("B6").Value = 1.542
$ws.Range("C6").Value = 0.005
$ws.Range("D6").Value = 0.6860000000000001
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = "Muhammad Iqbal Baqi"
$ws.Range("G6").Value = "Benar"

# Row 7
$ws.Range("A7").Value = "MIB_2.png"
$ws.Range("B7").Value = 1.204
$ws.Range("C7").Value = 0.004
$ws.Range("D7").Value = 0.759
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "Muhammad Iqbal Baqi"
$ws.Range("G7").Value = "Benar"

# Row 8
$ws.Range("A8").Value = "MIB_3.png"
$ws.Range("B8").Value = 1.702
$ws.Range("C8").Value = 0.006
$ws.Range("D8").Value = 0.97
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "Muhammad Iqbal Baqi"
$ws.Range("G8").Value = "Benar"

# Row 9
$ws.Range("A9").Value = "MIB_4.png"
$ws.Range("B9").Value = 1.205
$ws.Range("C9").Value = 0.004
$ws.Range("D9").Value = 0.785
$ws.Range("E9").Value = 0.714
$ws.Range("F9").Value = "Muhammad Iqbal Baqi"
$ws.Range("G9").Value = "Benar"

# Row 10
$ws.Range("A10").Value = "MIB_5.png"
$ws.Range("B10").Value = 1.5
$ws.Range("C10").Value = 0.005
$ws.Range("D10").Value = 0.733
$ws.Range("E10").Value = 0.714
$ws.Range("F10").Value = "Muhammad Iqbal Baqi"
$ws.Range("G10").Value = "Benar"

# Row 11
$ws.Range("A11").Value = "AAH_1.png"
$ws.Range("B11").Value = 0.856
$ws.Range("C11").Value = 0.003
$ws.Range("D11").Value = 0.794
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = "Andrea Ayunove Hutami"
$ws.Range("G11").Value = "Benar"

# Row 12
$ws.Range("A12").Value = "AAH_2.png"
$ws.Range("B12").Value = 1.114
$ws.Range("C12").Value = 0.004
$ws.Range("D12").Value = 0.946
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "Andrea Ayunove Hutami"
$ws.Range("G12").Value = "Benar"

# Row 13
$ws.Range("A13").Value = "AAH_3.png"
$ws.Range("B13").Value = 0.886
$ws.Range("C13").Value = 0.003
$ws.Range("D13").Value = 0.805
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "Andrea Ayunove Hutami"
$ws.Range("G13").Value = "Benar"

# Row 14
$ws.Range("A14").Value = "TI_1.png"
$ws.Range("B14").Value = 0.925
$ws.Range("C14").Value = 0.003
$ws.Range("D14").Value = 0.8100000000000001
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "Toni Ismail"
$ws.Range("G14").Value = "Benar"

# Row 15
$ws.Range("A15").Value = "TI_2.png"
$ws.Range("B15").Value = 0.9429999999999999
$ws.Range("C15").Value = 0.003
$ws.Range("D15").Value = 0.858
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = "Toni Ismail"
$ws.Range("G15").Value = "Benar"

# Row 16
$ws.Range("A16").Value = "TI_3.png"
$ws.Range("B16").Value = 0.727
$ws.Range("C16").Value = 0.002
$ws.Range("D16").Value = 0.948
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = "Toni Ismail"
$ws.Range("G16").Value = "Benar"

# Row 17
$ws.Range("A17").Value = "TI_4.png"
$ws.Range("B17").Value = 0.702
$ws.Range("C17").Value = 0.002
$ws.Range("D17").Value = 0.89
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = "Toni Ismail"
$ws.Range("G17").Value = "Benar"

# Row 18
$ws.Range("A18").Value = "TI_5.png"
$ws.Range("B18").Value = 1.07
$ws.Range("C18").Value = 0.004
$ws.Range("D18").Value = 0.91
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = "Toni Ismail"
$ws.Range("G18").Value = "Benar"

# Row 19
$ws.Range("A19").Value = "RAS_1.png"
$ws.Range("B19").Value = 0.634
$ws.Range("C19").Value = 0.002
$ws.Range("D19").Value = 0.844
$ws.Range("E19").Value = 0.857
$ws.Range("F19").Value = "Ridha Ayu Salsabila"
$ws.Range("G19").Value = "Benar"

# Row 20
$ws.Range("A20").Value = "RAS_2.png"
$ws.Range("B20").Value = 1.063
$ws.Range("C20").Value = 0.004
$ws.Range("D20").Value = 0.8139999999999999
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = "Ridha Ayu Salsabila"
$ws.Range("G20").Value = "Benar"

# Row 21
$ws.Range("A21").Value = "RAS_3.png"
$ws.Range("B21").Value = 0.519
$ws.Range("C21").Value = 0.002
$ws.Range("D21").Value = 0.798
$ws.Range("E21").Value = 0.857
$ws.Range("F21").Value = "Ridha Ayu Salsabila"
$ws.Range("G21").Value = "Benar"

# Row 22
$ws.Range("A22").Value = "RAS_4.png"
$ws.Range("B22").Value = 1.246
$ws.Range("C22").Value = 0.004
$ws.Range("D22").Value = 0.726
$ws.Range("E22").Value = 0.571
$ws.Range("F22").Value = "Ridha Ayu Salsabila"
$ws.Range("G22").Value = "Benar"

# Row 23
$ws.Range("A23").Value = "RAS_5.png"
$ws.Range("B23").Value = 1.179
$ws.Range("C23").Value = 0.004
$ws.Range("D23").Value = 0.847
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = "Ridha Ayu Salsabila"
$ws.Range("G23").Value = "Benar"

# Row 24
$ws.Range("A24").Value = "RR_1.png"
$ws.Range("B24").Value = 1.421
$ws.Range("C24").Value = 0.005
$ws.Range("D24").Value = 0.902
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = "Rafiqo Rapitasari"
$ws.Range("G24").Value = "Benar"

# Row 25
$ws.Range("A25").Value = "RR_2.png"
$ws.Range("B25").Value = 1.392
$ws.Range("C25").Value = 0.005
$ws.Range("D25").Value = 0.908
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = "Rafiqo Rapitasari"
$ws.Range("G25").Value = "Benar"

# Row 26
$ws.Range("A26").Value = "RR_3.png"
$ws.Range("B26").Value = 1.119
$ws.Range("C26").Value = 0.004
$ws.Range("D26").Value = 0.547
$ws.Range("E26").Value = 0.571
$ws.Range("F26").Value = "Rafiqo Rapitasari"
$ws.Range("G26").Value = "Benar"

# Row 27
$ws.Range("A27").Value = "RR_4.png"
$ws.Range("B27").Value = 1.374
$ws.Range("C27").Value = 0.005
$ws.Range("D27").Value = 0.905
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = "Rafiqo Rapitasari"
$ws.Range("G27").Value = "Benar"

# Row 28
$ws.Range("A28").Value = "RR_5.png"
$ws.Range("B28").Value = 1.457
$ws.Range("C28").Value = 0.005
$ws.Range("D28").Value = 0.904
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = "Rafiqo Rapitasari"
$ws.Range("G28").Value = "Benar"

# Row 29
$ws.Range("A29").Value = "AR_1.png"
$ws.Range("B29").Value = 0.863
$ws.Range("C29").Value = 0.003
$ws.Range("D29").Value = 0.89
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = "Arizli Romadhon"
$ws.Range("G29").Value = "Benar"

# Row 30
$ws.Range("A30").Value = "GA_1.png"
$ws.Range("B30").Value = 1.58
$ws.Range("C30").Value = 0.005
$ws.Range("D30").Value = 0.898
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = "Gege Ardiyansyah"
$ws.Range("G30").Value = "Benar"

# Row 31
$ws.Range("A31").Value = "GA_2.png"
$ws.Range("B31").Value = 0.6879999999999999
$ws.Range("C31").Value = 0.002
$ws.Range("D31").Value = 0.832
$ws.Range("E31").Value = 0.286
$ws.Range("F31").Value = "Tidak Diketahui"
$ws.Range("G31").Value = "Salah"

# Row 32
$ws.Range("A32").Value = "GA_3.png"
$ws.Range("B32").Value = 0.727
$ws.Range("C32").Value = 0.002
$ws.Range("D32").Value = 0.79
$ws.Range("E32").Value = 0.286
$ws.Range("F32").Value = "Tidak Diketahui"
$ws.Range("G32").Value = "Salah"

# Row 33
$ws.Range("A33").Value = "FY_1.png"
$ws.Range("B33").Value = 1.09
$ws.Range("C33").Value = 0.004
$ws.Range("D33").Value = 0.888
$ws.Range("E33").Value = 0.571
$ws.Range("F33").Value = "Fanny Yusuf"
$ws.Range("G33").Value = "Benar"

# Row 34
$ws.Range("A34").Value = "FY_2.png"
$ws.Range("B34").Value = 1.602
$ws.Range("C34").Value = 0.005
$ws.Range("D34").Value = 0.872
$ws.Range("E34").Value = 0.857
$ws.Range("F34").Value = "Fanny Yusuf"
$ws.Range("G34").Value = "Benar"

# Row 35
$ws.Range("A35").Value = "FY_3.png"
$ws.Range("B35").Value = 1.477
$ws.Range("C35").Value = 0.005
$ws.Range("D35").Value = 0.893
$ws.Range("E35").Value = 0.857
$ws.Range("F35").Value = "Fanny Yusuf"
$ws.Range("G35").Value = "Benar"

# Row 36
$ws.Range("A36").Value = "FY_4.png"
$ws.Range("B36").Value = 1.257
$ws.Range("C36").Value = 0.004
$ws.Range("D36").Value = 0.834
$ws.Range("E36").Value = 0.714
$ws.Range("F36").Value = "Fanny Yusuf"
$ws.Range("G36").Value = "Benar"

# Row 37
$ws.Range("A37").Value = "TO_1.png"
$ws.Range("B37").Value = 0.781
$ws.Range("C37").Value = 0.003
$ws.Range("D37").Value = 0.806
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = "Tiara Oktavian"
$ws.Range("G37").Value = "Benar"

# Row 38
$ws.Range("A38").Value = "TO_2.png"
$ws.Range("B38").Value = 1.068
$ws.Range("C38").Value = 0.004
$ws.Range("D38").Value = 0.857
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = "Tiara Oktavian"
$ws.Range("G38").Value = "Benar"

# Row 39
$ws.Range("A39").Value = "TO_3.png"
$ws.Range("B39").Value = 0.843
$ws.Range("C39").Value = 0.003
$ws.Range("D39").Value = 0.858
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = "Tiara Oktavian"
$ws.Range("G39").Value = "Benar"

# Row 40
$ws.Range("A40").Value = "TO_4.png"
$ws.Range("B40").Value = 3.224
$ws.Range("C40").Value = 0.01
$ws.Range("D40").Value = 0.514
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = "Tiara Oktavian"
$ws.Range("G40").Value = "Benar"

# Row 41
$ws.Range("A41").Value = "TO_5.png"
$ws.Range("B41").Value = 3.124
$ws.Range("C41").Value = 0.01
$ws.Range("D41").Value = 0.503
$ws.Range("E41").Value = 0.714
$ws.Range("F41").Value = "Tiara Oktavian"
$ws.Range("G41").Value = "Benar"

# Row 42
$ws.Range("A42").Value = "TD_1.png"
$ws.Range("B42").Value = 2.34
$ws.Range("C42").Value = 0.008
$ws.Range("D42").Value = 0.433
$ws.Range("E42").Value = 0.429
$ws.Range("F42").Value = "Tidak Diketahui"
$ws.Range("G42").Value = "Benar"

# Row 43
$ws.Range("A43").Value = "TD_2.png"
$ws.Range("B43").Value = 2.479
$ws.Range("C43").Value = 0.008
$ws.Range("D43").Value = 0.462
$ws.Range("E43").Value = 0.429
$ws.Range("F43").Value = "Tidak Diketahui"
$ws.Range("G43").Value = "Benar"

# Row 44
$ws.Range("A44").Value = "TD_3.png"
$ws.Range("B44").Value = 1.081
$ws.Range("C44").Value = 0.004
$ws.Range("D44").Value = 0.732
$ws.Range("E44").Value = 0.286
$ws.Range("F44").Value = "Tidak Diketahui"
$ws.Range("G44").Value = "Benar"

# Row 45
$ws.Range("A45").Value = "TD_4.png"
$ws.Range("B45").Value = 1.008
$ws.Range("C45").Value = 0.003
$ws.Range("D45").Value = 0.718
$ws.Range("E45").Value = 0.286
$ws.Range("F45").Value = "Tidak Diketahui"
$ws.Range("G45").Value = "Benar"
